# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Source cells store prices as *text* (not numbers) -- e.g. "1.599.14" uses a dot
# as a thousands separator, and values like "15.30"/"1.50" need the trailing zero
# kept literally. Plain Excel autodetects the latter kind ("19.53", "211.21", ...)
# as numbers, so those are written via .Formula with a leading apostrophe to force
# text storage, matching the original inlineStr cells exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.711.04'
$ws.Range('E2').Value = '  +0.31%  '

$ws.Range('D3').Value = '1.598.59'
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').Formula = '''211.21'
$ws.Range('E5').Value = '  +0.03%  '

$ws.Range('E6').Value = '  -0.66%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('E8').Value = '  +0.34%  '

$ws.Range('E9').Value = '  +1.06%  '

$ws.Range('D10').Formula = '''19.53'
$ws.Range('E10').Value = '  +0.85%  '

$ws.Range('E11').Value = '  +0.50%  '

$ws.Range('D12').Value = '1.823.71'
$ws.Range('E12').Value = '  +0.24%  '

$ws.Range('D13').Value = '1.615.63'
$ws.Range('E13').Value = '  +1.50%  '

$ws.Range('E14').Value = '  +0.54%  '

$ws.Range('E15').Value = '  +0.42%  '

$ws.Range('D16').Formula = '''65.28'
$ws.Range('E16').Value = '  +0.92%  '

$ws.Range('D17').Value = '26.683.80'
$ws.Range('E17').Value = '  +0.33%  '

$ws.Range('D18').Value = '0.0₃0760'
$ws.Range('E18').Value = '  +4.12%  '

$ws.Range('D19').Formula = '''209.87'
$ws.Range('E19').Value = '  +0.92%  '

$ws.Range('D21').Formula = '''7.13'
$ws.Range('E21').Value = '  +3.29%  '

$ws.Range('E22').Value = '  +0.61%  '

$ws.Range('D23').Formula = '''2.29'
$ws.Range('E23').Value = '  -0.32%  '

$ws.Range('E24').Value = '  +0.90%  '

$ws.Range('D25').Formula = '''143.12'
$ws.Range('E25').Value = '  -1.58%  '

$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('E27').Value = '  -0.21%  '

$ws.Range('E28').Value = '  +0.32%  '

$ws.Range('D29').Formula = '''15.30'
$ws.Range('E29').Value = '  +0.12%  '

$ws.Range('D30').Formula = '''0.0518'
$ws.Range('E30').Value = '  +2.83%  '

$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('E32').Value = '  +0.40%  '

$ws.Range('E33').Value = '  +1.58%  '

$ws.Range('D34').Value = '1.288.43'
$ws.Range('E34').Value = '  +0.27%  '

$ws.Range('D35').Formula = '''0.618'
$ws.Range('E35').Value = '  -5.77%  '

$ws.Range('E36').Value = '  +1.14%  '

$ws.Range('D37').Formula = '''1.50'
$ws.Range('E37').Value = '  +0.20%  '

$ws.Range('E38').Value = '  -0.23%  '

$ws.Range('E39').Value = '  +17.61%  '

$ws.Range('D40').Formula = '''0.826'
$ws.Range('E40').Value = '  -1.52%  '

$ws.Range('E41').Value = '  -0.11%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Formula = '''2.19'
$ws.Range('E42').Value = '  -0.46%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Formula = '''0.783'
$ws.Range('E43').Value = '  -0.27%  '

$ws.Range('D44').Formula = '''63.05'
$ws.Range('E44').Value = '  -0.98%  '

$ws.Range('D45').Value = '1.736.60'
$ws.Range('E45').Value = '  +0.28%  '

$ws.Range('D46').Formula = '''91.25'
$ws.Range('E46').Value = '  +1.75%  '

$ws.Range('E47').Value = '  -1.55%  '

$ws.Range('E48').Value = '  -0.89%  '

$ws.Range('E49').Value = '  +0.51%  '

$ws.Range('E50').Value = '  +0.14%  '

$ws.Range('D51').Formula = '''7.36'
$ws.Range('E51').Value = '  -1.22%  '
